# Lattice multiplication exercises table update
# Replaces the 15 practice problems (3 cols x 5 rows) with new operands
# and regenerates the corresponding lattice grid labels in each cell.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$br = [char]11   # vertical-tab char Word uses for a <w:br/> line break

# Row 1, Col 1: "87 x 99" -> "51 x 21"
$cell = $t.Cell(1, 1)
$para = $cell.Range.Paragraphs.Item(1)
$para.Range.Text = "51 x 21" + $br + "  2    1" + $br + "  ----" + $br + "5|    |" + $br + "1|    |"

# Row 1, Col 2: "80 x 50" -> "94 x 98"
$cell = $t.Cell(1, 2)
$para = $cell.Range.Paragraphs.Item(1)
$para.Range.Text = "94 x 98" + $br + "  9    8" + $br + "  ----" + $br + "9|    |" + $br + "4|    |"

# Row 1, Col 3: "20 x 97" -> "42 x 89"
$cell = $t.Cell(1, 3)
$para = $cell.Range.Paragraphs.Item(1)
$para.Range.Text = "42 x 89" + $br + "  8    9" + $br + "  ----" + $br + "4|    |" + $br + "2|    |"

# Row 2, Col 1: "41 x 43" -> "64 x 88"
$cell = $t.Cell(2, 1)
$para = $cell.Range.Paragraphs.Item(1)
$para.Range.Text = "64 x 88" + $br + "  8    8" + $br + "  ----" + $br + "6|    |" + $br + "4|    |"

# Row 2, Col 2: "34 x 22" -> "62 x 84"
$cell = $t.Cell(2, 2)
$para = $cell.Range.Paragraphs.Item(1)
$para.Range.Text = "62 x 84" + $br + "  8    4" + $br + "  ----" + $br + "6|    |" + $br + "2|    |"

# Row 2, Col 3: "54 x 83" -> "79 x 19"
$cell = $t.Cell(2, 3)
$para = $cell.Range.Paragraphs.Item(1)
$para.Range.Text = "79 x 19" + $br + "  1    9" + $br + "  ----" + $br + "7|    |" + $br + "9|    |"

# Row 3, Col 1: "60 x 90" -> "94 x 44"
$cell = $t.Cell(3, 1)
$para = $cell.Range.Paragraphs.Item(1)
$para.Range.Text = "94 x 44" + $br + "  4    4" + $br + "  ----" + $br + "9|    |" + $br + "4|    |"

# Row 3, Col 2: "34 x 82" -> "87 x 95"
$cell = $t.Cell(3, 2)
$para = $cell.Range.Paragraphs.Item(1)
$para.Range.Text = "87 x 95" + $br + "  9    5" + $br + "  ----" + $br + "8|    |" + $br + "7|    |"

# Row 3, Col 3: "59 x 78" -> "63 x 22"
$cell = $t.Cell(3, 3)
$para = $cell.Range.Paragraphs.Item(1)
$para.Range.Text = "63 x 22" + $br + "  2    2" + $br + "  ----" + $br + "6|    |" + $br + "3|    |"

# Row 4, Col 1: "62 x 66" -> "40 x 49"
$cell = $t.Cell(4, 1)
$para = $cell.Range.Paragraphs.Item(1)
$para.Range.Text = "40 x 49" + $br + "  4    9" + $br + "  ----" + $br + "4|    |" + $br + "0|    |"

# Row 4, Col 2: "16 x 91" -> "31 x 31"
$cell = $t.Cell(4, 2)
$para = $cell.Range.Paragraphs.Item(1)
$para.Range.Text = "31 x 31" + $br + "  3    1" + $br + "  ----" + $br + "3|    |" + $br + "1|    |"

# Row 4, Col 3: "25 x 11" -> "30 x 26"
$cell = $t.Cell(4, 3)
$para = $cell.Range.Paragraphs.Item(1)
$para.Range.Text = "30 x 26" + $br + "  2    6" + $br + "  ----" + $br + "3|    |" + $br + "0|    |"

# Row 5, Col 1: "38 x 72" -> "90 x 44"
$cell = $t.Cell(5, 1)
$para = $cell.Range.Paragraphs.Item(1)
$para.Range.Text = "90 x 44" + $br + "  4    4" + $br + "  ----" + $br + "9|    |" + $br + "0|    |"

# Row 5, Col 2: "21 x 58" -> "76 x 81"
$cell = $t.Cell(5, 2)
$para = $cell.Range.Paragraphs.Item(1)
$para.Range.Text = "76 x 81" + $br + "  8    1" + $br + "  ----" + $br + "7|    |" + $br + "6|    |"

# Row 5, Col 3: "79 x 49" -> "35 x 29"
$cell = $t.Cell(5, 3)
$para = $cell.Range.Paragraphs.Item(1)
$para.Range.Text = "35 x 29" + $br + "  2    9" + $br + "  ----" + $br + "3|    |" + $br + "5|    |"

Write-Output "Done"
